$d = $word.ActiveDocument

# The original single paragraph ("Hello GIT") carries the hidden
# "_GoBack" bookmark around its whole text. The edit splits that
# paragraph into three paragraphs (same centered/bold/72pt formatting)
# and re-homes the (now empty) "_GoBack" bookmark at the very end of
# the third/last paragraph.

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p1 = $d.Paragraphs(1)

# Split after "Hello GIT" -> new paragraph inherits pPr/rPr formatting.
$p1.Range.InsertParagraphAfter()
$d.Paragraphs(2).Range.Text = ";)))))))"

# Split again -> third paragraph.
$d.Paragraphs(2).Range.InsertParagraphAfter()
# Temporary trailing marker char "Z" lets us park a zero-length
# bookmark mid-paragraph (placing it exactly at the paragraph-mark
# boundary is unreliable), then we trim the marker back off while the
# bookmark (anchored strictly before it) stays put.
$d.Paragraphs(3).Range.Text = "))Z"

$p3 = $d.Paragraphs(3)
$markerStart = $p3.Range.End - 2   # position just before "Z", after "))"
$bmRange = $d.Range($markerStart, $markerStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

$delRange = $d.Range($markerStart, $markerStart + 1)
$delRange.Delete()
